# "Added 1.1.0 of term"
# The workbook's "Metadata" sheet is a two-column Property/Value table.
# Bump the Version value and refresh the Date value that goes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$versionCell = $ws.Range("A1:A14").Find("Version").Offset(0, 1)
$versionCell.Value = "1.1.0"

$dateCell = $ws.Range("A1:A14").Find("Date").Offset(0, 1)
$dateCell.Value = "2023-07-10T23:08:03+02:00"
